$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B18").Value = "Allan Cupertino-Máquinas Elétri"
$ws.Range("C18").Value = "[Emerson-Eletrônica Básica, Allan Cupertino-Instalções Elétricas]"
$ws.Range("D18").Value = "[Cleidson-Automação Industrial, Guilherme-Eletrohidráulica, -, Cláudio-Tecnologia da Soldagem]"
$ws.Range("E18").Value = "[Emerson-Eletrônica Básica, Weslei-CAD]"
$ws.Range("F18").Value = "[Paulo Rob.-CAM, Leandro-Sistemas de Refrigeração, Guilherme-Eletrohidráulica, -]"

$ws.Range("B19").Value = "Allan Cupertino-Máquinas Elétri"
$ws.Range("C19").Value = "[Emerson-Eletrônica Básica, Allan Cupertino-Instalções Elétricas]"
$ws.Range("D19").Value = "[Cleidson-Automação Industrial, Guilherme-Eletropneumática, Cláudio-Tecnologia da Soldagem, Paulo Rob.-CAM]"
$ws.Range("E19").Value = "[Allan Cupertino-Instalções Elétricas, Weslei-CAD]"
$ws.Range("F19").Value = "[Leandro-Sistemas de Refrigeração, Paulo Rob.-CAM, Guilherme-Eletrohidráulica, -]"

$ws.Range("B20").Value = "Andre B.-Circuitos Elétrico"
$ws.Range("C20").Value = "[Allan Cupertino-Lab. De Máquinas elétricas, João Paulo-Lab. Circuitos Elétricos]"
$ws.Range("D20").Value = "[Cleidson-Automação Industrial, Guilherme-Eletropneumática, -, -]"
$ws.Range("E20").Value = "[Allan Cupertino-Instalções Elétricas, Weslei-CAD]"
$ws.Range("F20").Value = "[Guilherme-Eletropneumática, Paulo Rob.-CAM, -, -]"

$ws.Range("B21").Value = "Andre B.-Circuitos Elétrico"
$ws.Range("C21").Value = "[Allan Cupertino-Lab. De Máquinas elétricas, João Paulo-Lab. Circuitos Elétricos]"
$ws.Range("D21").Value = "[Cleidson-Automação Industrial, Guilherme-Eletropneumática, Cláudio-Tecnologia da Soldagem, Leandro-Sistemas de Refrigeração]"
$ws.Range("E21").Value = "[Emerson-Eletrônica Básica, Weslei-CAD]"
$ws.Range("F21").Value = "[-, Leandro-Sistemas de Refrigeração, Cláudio-Tecnologia da Soldagem, Guilherme-Eletrohidráulica]"
